# last_week_prices.xlsx — "first kinda working algo"
# Column A used to hold actual dates (rows 4,6,...,26) with a date-number
# format (style index 1, numFmtId 16 "d-mmm"), while the header-ish rows
# (2, 28..48) held hour-range labels as plain text with no explicit style.
# The edit turns the whole column into hour-range labels ("0 - 1", "1 - 2",
# ... "23 - 24") and repurposes style 1 as a text format (numFmtId 49 "@")
# applied uniformly down column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-point the existing "date" style (index 1) at a text number format
# instead of a date format, since the column no longer holds real dates.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A24").NumberFormat = "@"
$ws.Range("A26").NumberFormat = "@"
$ws.Range("A28").NumberFormat = "@"
$ws.Range("A30").NumberFormat = "@"
$ws.Range("A32").NumberFormat = "@"
$ws.Range("A34").NumberFormat = "@"
$ws.Range("A36").NumberFormat = "@"
$ws.Range("A38").NumberFormat = "@"
$ws.Range("A40").NumberFormat = "@"
$ws.Range("A42").NumberFormat = "@"
$ws.Range("A44").NumberFormat = "@"
$ws.Range("A46").NumberFormat = "@"
$ws.Range("A48").NumberFormat = "@"

# Replace the date values in column A with hour-range labels, typed in the
# same (slightly out-of-order) sequence as the original editing session.
$ws.Range("A4").Value = "1 - 2"
$ws.Range("A8").Value = "3 - 4"
$ws.Range("A10").Value = "4 - 5"
$ws.Range("A6").Value = "2 - 3"
$ws.Range("A14").Value = "6 - 7"
$ws.Range("A16").Value = "7 - 8"
$ws.Range("A20").Value = "9 - 10"
$ws.Range("A24").Value = "11 - 12"
$ws.Range("A26").Value = "12 - 13"
$ws.Range("A12").Value = "5 - 6"
$ws.Range("A18").Value = "8 - 9"
$ws.Range("A22").Value = "10 - 11"

# Give column A a wider, fixed width so the new labels read comfortably.
$ws.Columns.Item(1).ColumnWidth = 52.3

# Print setup: portrait orientation.
$ws.PageSetup.Orientation = 1

# Leave the selection where the editing session ended up.
$ws.Range("A23").Select()
